$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.614.47"
$ws.Range("E2").Value = "  -0.59%  "
$ws.Range("D3").Value = "2.721.44"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.51%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -1.47%  "
$ws.Range("E9").Value = "  -2.26%  "
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.58"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.373"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.96%  "
$ws.Range("D13").Value = "3.204.49"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.64"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "63.483.82"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("D17").Value = "2.724.43"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.18"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  -3.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "348.31"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.73%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  -2.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.09%  "
$ws.Range("D28").Value = "0.0₃0884"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("E29").Value = "  +9.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.84"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.86"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.17%  "
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.85"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.07%  "
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("E37").Value = "  +0.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "347.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.958"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.16%  "
$ws.Range("E40").Value = "  -1.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0573"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.627"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "132.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.98%  "
$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.07"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0246"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.86%  "
$ws.Range("E51").Value = "  -2.82%  "
